# Add data for 2024-09-01
# Updates the 2024 (column K) running-total figures across all affected worksheets
# to incorporate September 2024 crime counts.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Cells.Item(2, 11).Value = 5355
$ws.Cells.Item(3, 11).Value = 5511
$ws.Cells.Item(4, 11).Value = 1150
$ws.Cells.Item(5, 11).Value = 395
$ws.Cells.Item(6, 11).Value = 6135
$ws.Cells.Item(7, 11).Value = 18546

$ws = $wb.Worksheets.Item("Austin")
$ws.Cells.Item(2, 11).Value = 341
$ws.Cells.Item(7, 11).Value = 1237

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Cells.Item(2, 11).Value = 140
$ws.Cells.Item(5, 11).Value = 10
$ws.Cells.Item(7, 11).Value = 413

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Cells.Item(2, 11).Value = 218
$ws.Cells.Item(3, 11).Value = 293
$ws.Cells.Item(4, 11).Value = 36
$ws.Cells.Item(6, 11).Value = 234
$ws.Cells.Item(7, 11).Value = 798

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Cells.Item(2, 11).Value = 180
$ws.Cells.Item(3, 11).Value = 206
$ws.Cells.Item(6, 11).Value = 182
$ws.Cells.Item(7, 11).Value = 626

$ws = $wb.Worksheets.Item("New City")
$ws.Cells.Item(2, 11).Value = 135
$ws.Cells.Item(4, 11).Value = 14
$ws.Cells.Item(7, 11).Value = 421

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Cells.Item(6, 11).Value = 25
$ws.Cells.Item(7, 11).Value = 73

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Cells.Item(2, 11).Value = 163
$ws.Cells.Item(4, 11).Value = 70
$ws.Cells.Item(6, 11).Value = 134
$ws.Cells.Item(7, 11).Value = 550
$ws.Cells.Item(8, 11).Value = 1237
$ws.Cells.Item(18, 11).Value = 125
$ws.Cells.Item(19, 11).Value = 546
$ws.Cells.Item(20, 11).Value = 433
$ws.Cells.Item(21, 11).Value = 59
$ws.Cells.Item(29, 11).Value = 994
$ws.Cells.Item(30, 11).Value = 73
$ws.Cells.Item(31, 11).Value = 201
$ws.Cells.Item(33, 11).Value = 798
$ws.Cells.Item(37, 11).Value = 626
$ws.Cells.Item(41, 11).Value = 129
$ws.Cells.Item(42, 11).Value = 687
$ws.Cells.Item(47, 11).Value = 126
$ws.Cells.Item(52, 11).Value = 486
$ws.Cells.Item(54, 11).Value = 362
$ws.Cells.Item(55, 11).Value = 206
$ws.Cells.Item(57, 11).Value = 71
$ws.Cells.Item(59, 11).Value = 30
$ws.Cells.Item(63, 11).Value = 51
$ws.Cells.Item(64, 11).Value = 119
$ws.Cells.Item(65, 11).Value = 421
$ws.Cells.Item(67, 11).Value = 706
$ws.Cells.Item(72, 11).Value = 89
$ws.Cells.Item(76, 11).Value = 257
$ws.Cells.Item(79, 11).Value = 460
$ws.Cells.Item(83, 11).Value = 413
$ws.Cells.Item(84, 11).Value = 140
$ws.Cells.Item(85, 11).Value = 874
$ws.Cells.Item(86, 11).Value = 124
$ws.Cells.Item(88, 11).Value = 203
$ws.Cells.Item(89, 11).Value = 271
$ws.Cells.Item(91, 11).Value = 206
$ws.Cells.Item(94, 11).Value = 249
$ws.Cells.Item(96, 11).Value = 201
$ws.Cells.Item(97, 11).Value = 148
$ws.Cells.Item(100, 11).Value = 35
$ws.Cells.Item(101, 11).Value = 18546

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Cells.Item(2, 11).Value = 67
$ws.Cells.Item(7, 11).Value = 201

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Cells.Item(3, 11).Value = 253
$ws.Cells.Item(5, 11).Value = 16
$ws.Cells.Item(6, 11).Value = 197
$ws.Cells.Item(7, 11).Value = 706

$ws = $wb.Worksheets.Item("South Deering")
$ws.Cells.Item(2, 11).Value = 45
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(7, 11).Value = 140

$ws = $wb.Worksheets.Item("Loop")
$ws.Cells.Item(6, 11).Value = 193
$ws.Cells.Item(7, 11).Value = 362

$ws = $wb.Worksheets.Item("Englewood")
$ws.Cells.Item(2, 11).Value = 283
$ws.Cells.Item(3, 11).Value = 357
$ws.Cells.Item(6, 11).Value = 279
$ws.Cells.Item(7, 11).Value = 994

$ws = $wb.Worksheets.Item("Chatham")
$ws.Cells.Item(2, 11).Value = 162
$ws.Cells.Item(7, 11).Value = 546

$ws = $wb.Worksheets.Item("River North")
$ws.Cells.Item(3, 11).Value = 47
$ws.Cells.Item(7, 11).Value = 257

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Cells.Item(3, 11).Value = 38
$ws.Cells.Item(7, 11).Value = 134

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Cells.Item(6, 11).Value = 50
$ws.Cells.Item(7, 11).Value = 129

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Cells.Item(3, 11).Value = 213
$ws.Cells.Item(6, 11).Value = 261
$ws.Cells.Item(7, 11).Value = 687

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Cells.Item(3, 11).Value = 57
$ws.Cells.Item(4, 11).Value = 10
$ws.Cells.Item(6, 11).Value = 71
$ws.Cells.Item(7, 11).Value = 206

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Cells.Item(6, 11).Value = 87
$ws.Cells.Item(7, 11).Value = 201

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Cells.Item(2, 11).Value = 52
$ws.Cells.Item(3, 11).Value = 97
$ws.Cells.Item(4, 11).Value = 7
$ws.Cells.Item(7, 11).Value = 206

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Cells.Item(6, 11).Value = 36
$ws.Cells.Item(7, 11).Value = 59

$ws = $wb.Worksheets.Item("Roseland")
$ws.Cells.Item(2, 11).Value = 155
$ws.Cells.Item(6, 11).Value = 112
$ws.Cells.Item(7, 11).Value = 460

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Cells.Item(2, 11).Value = 23
$ws.Cells.Item(7, 11).Value = 119

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Cells.Item(2, 11).Value = 144
$ws.Cells.Item(3, 11).Value = 137
$ws.Cells.Item(6, 11).Value = 125
$ws.Cells.Item(7, 11).Value = 433

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Cells.Item(2, 11).Value = 34
$ws.Cells.Item(7, 11).Value = 125

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Cells.Item(3, 11).Value = 5
$ws.Cells.Item(7, 11).Value = 35

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Cells.Item(3, 11).Value = 176
$ws.Cells.Item(6, 11).Value = 147
$ws.Cells.Item(7, 11).Value = 550

$ws = $wb.Worksheets.Item("West Loop")
$ws.Cells.Item(2, 11).Value = 70
$ws.Cells.Item(6, 11).Value = 110
$ws.Cells.Item(7, 11).Value = 249

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Cells.Item(2, 11).Value = 36
$ws.Cells.Item(7, 11).Value = 126

$ws = $wb.Worksheets.Item("Montclare")
$ws.Cells.Item(6, 11).Value = 9
$ws.Cells.Item(7, 11).Value = 30

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Cells.Item(3, 11).Value = 38
$ws.Cells.Item(7, 11).Value = 163

$ws = $wb.Worksheets.Item("West Town")
$ws.Cells.Item(3, 11).Value = 29
$ws.Cells.Item(7, 11).Value = 148

$ws = $wb.Worksheets.Item("United Center")
$ws.Cells.Item(2, 11).Value = 53
$ws.Cells.Item(3, 11).Value = 61
$ws.Cells.Item(7, 11).Value = 203

$ws = $wb.Worksheets.Item("Uptown")
$ws.Cells.Item(6, 11).Value = 81
$ws.Cells.Item(7, 11).Value = 271

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Cells.Item(6, 11).Value = 31
$ws.Cells.Item(7, 11).Value = 124

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Cells.Item(4, 11).Value = 5
$ws.Cells.Item(6, 11).Value = 32
$ws.Cells.Item(7, 11).Value = 71

$ws = $wb.Worksheets.Item("South Shore")
$ws.Cells.Item(2, 11).Value = 292
$ws.Cells.Item(3, 11).Value = 295
$ws.Cells.Item(7, 11).Value = 874

$ws = $wb.Worksheets.Item("Old Town")
$ws.Cells.Item(6, 11).Value = 46
$ws.Cells.Item(7, 11).Value = 89

$ws = $wb.Worksheets.Item("Little Village")
$ws.Cells.Item(2, 11).Value = 131
$ws.Cells.Item(6, 11).Value = 177
$ws.Cells.Item(7, 11).Value = 486

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Cells.Item(2, 11).Value = 24
$ws.Cells.Item(7, 11).Value = 70

